$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update id (B) and speaker_variant (C) columns for rows 2-36
$ws.Range("B2").Value = '#gisippus'
$ws.Range("C2").Value = 'Gisippus'
$ws.Range("B3").Value = '#spaernoy'
$ws.Range("C3").Value = 'Spaernoy'
$ws.Range("B4").Value = '#m.-varro'
$ws.Range("C4").Value = 'M. Varro'
$ws.Range("B5").Value = '#seld'
$ws.Range("C5").Value = 'Seld'
$ws.Range("B6").Value = '#dienaer'
$ws.Range("C6").Value = 'Dienaer'
$ws.Range("B7").Value = '#erato'
$ws.Range("C7").Value = 'Erato'
$ws.Range("B8").Value = '#leeghwaghen'
$ws.Range("C8").Value = 'Leeghwaghen'
$ws.Range("B9").Value = '#spaer-noy'
$ws.Range("C9").Value = 'Spaer-noy'
$ws.Range("B10").Value = '#dien'
$ws.Range("C10").Value = 'Dien'
$ws.Range("B11").Value = '#spaer'
$ws.Range("C11").Value = 'Spaer'
$ws.Range("B12").Value = '#octavianus'
$ws.Range("C12").Value = 'Octavianus'
$ws.Range("B13").Value = '#varro'
$ws.Range("C13").Value = 'Varro'
$ws.Range("B14").Value = '#p.-ambrosius'
$ws.Range("C14").Value = 'P. Ambrosius'
$ws.Range("B15").Value = '#schryver'
$ws.Range("C15").Value = 'Schryver'
$ws.Range("B16").Value = '#selden-thuijs'
$ws.Range("C16").Value = 'Selden-thuijs'
$ws.Range("B17").Value = '#pub.-ambrosius'
$ws.Range("C17").Value = 'Pub. Ambrosius'
$ws.Range("B18").Value = '#gysippus'
$ws.Range("C18").Value = 'Gysippus'
$ws.Range("B19").Value = '#spaer-noy'
$ws.Range("C19").Value = 'Spaer noy'
$ws.Range("B20").Value = '#dul-cop'
$ws.Range("C20").Value = 'Dul-cop'
$ws.Range("B21").Value = '#gis'
$ws.Range("C21").Value = 'Gis'
$ws.Range("B22").Value = '#schoutet'
$ws.Range("C22").Value = 'Schoutet'
$ws.Range("B23").Value = '#die'
$ws.Range("C23").Value = 'Die'
$ws.Range("B24").Value = '#titus'
$ws.Range("C24").Value = 'Titus'
$ws.Range("B25").Value = '#schou'
$ws.Range("C25").Value = 'Schou'
$ws.Range("B26").Value = '#sophronia'
$ws.Range("C26").Value = 'Sophronia'
$ws.Range("B27").Value = '#reysigher'
$ws.Range("C27").Value = 'Reysigher'
$ws.Range("B28").Value = '#gisip'
$ws.Range("C28").Value = 'Gisip'
$ws.Range("B29").Value = '#sophro'
$ws.Range("C29").Value = 'Sophro'
$ws.Range("B30").Value = '#pub.-amb'
$ws.Range("C30").Value = 'pub. Amb'
$ws.Range("B31").Value = '#thalia'
$ws.Range("C31").Value = 'Thalia'
$ws.Range("B32").Value = '#aristippus'
$ws.Range("C32").Value = 'Aristippus'
$ws.Range("B33").Value = '#tit'
$ws.Range("C33").Value = 'Tit'
$ws.Range("B34").Value = '#selden-thuys'
$ws.Range("C34").Value = 'Selden-thuys'
$ws.Range("B35").Value = '#selden-thuijs'
$ws.Range("C35").Value = 'Selden thuijs'
$ws.Range("B36").Value = '#verneem-al'
$ws.Range("C36").Value = 'Verneem-al'

# Clear is_prefered (D) column for rows 2-28 (no is_pref)
$ws.Range("D2:D28").ClearContents()

